$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header labels in F1:K1 to reflect "Prior Year" / "Curr. Year" naming
# instead of the hard-coded "2023" / "2024" years.
$ws.Range("F1").Value = "Prior Year LMV"
$ws.Range("G1").Value = "Prior Year BMV"
$ws.Range("H1").Value = "Prior Year Total"
$ws.Range("I1").Value = "Curr. Year LMV"
$ws.Range("J1").Value = "Curr. Year BMV"
$ws.Range("K1").Value = "Curr. Year Total"

# Declutter / fix column widths: give columns D and F:K explicit widths
# (no longer auto "best fit") now that the headers are longer.
$ws.Columns.Item(4).ColumnWidth = 9.5
$ws.Columns.Item(6).ColumnWidth = 14.833333333333334
$ws.Columns.Item(7).ColumnWidth = 15
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668
$ws.Columns.Item(9).ColumnWidth = 15.333333333333334
$ws.Columns.Item(10).ColumnWidth = 15.5
$ws.Columns.Item(11).ColumnWidth = 15.166666666666666

# Reset the view to open at A1 instead of scrolled to C1.
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
